# "further cleaning to metadata"
#
#   1. libraryProtocol value: "E7760" -> "E7420" for every data row (K2:K41).
#   2. roboticLibraryPrep column (L2:L41): replace the literal FALSE boolean
#      with the equivalent =FALSE() formula.
#   3. Update the saved selection from L2:L41 to K2:K41 (and scroll so row 9
#      is at the top of the view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 41

for ($row = 2; $row -le $lastRow; $row++) {
    # Column K = 11 -> libraryProtocol
    $ws.Cells.Item($row, 11).Value = "E7420"

    # Column L = 12 -> roboticLibraryPrep
    $ws.Cells.Item($row, 12).Formula = "=FALSE()"
}

# Selection / scroll position bookkeeping.
$ws.Range("K2:K41").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
